# Applies the edits described by the commit:
#  1. Update the cached datetimeFigureOut field text ("1/17/2017" -> "10/12/2016")
#     on the slide master and every slide layout's Date placeholder.
#  2. Fix capitalization in the "Many machine learning algorithms ..." sentence
#     on slide 3 ("Machine Learning Algorithms").
#  3. Merge the two "Feature" bullet points on slide 4 into a single bullet.

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16

function Set-DateFieldText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $isDatePlaceholder = $false
            if ($shape.Type -eq 14) {
                try {
                    if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                        $isDatePlaceholder = $true
                    }
                } catch {
                }
            }
            if ($isDatePlaceholder) {
                $shape.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# 1. Slide master date placeholder.
Set-DateFieldText $p.SlideMaster.Shapes "10/12/2016"

# 1. Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DateFieldText $layout.Shapes "10/12/2016"
}

# 2. Slide 3: "Many machine learning algorithms will take a similar approach:"
#    -> "Many machine Learning Algorithms will take a similar approach:"
$slide3 = $p.Slides.Item(3)
$contentShape3 = $slide3.Shapes.Item(2)
$run1 = $contentShape3.TextFrame.TextRange.Runs(1, 1)
$run1.Text = "Many machine Learning Algorithms will take a similar approach:"

# 3. Slide 4: merge the two "Feature"/"Features" bullets into one bullet.
$slide4 = $p.Slides.Item(4)
$contentShape4 = $slide4.Shapes.Item(2)
$tr4 = $contentShape4.TextFrame.TextRange
$featuresPara = $tr4.Paragraphs(5, 1)
$featuresPara.Delete()
$featurePara = $tr4.Paragraphs(4, 1)
$featureRun = $featurePara.Runs(1, 1)
$featureRun.Text = "Features – size of house, number of rooms, location"
